# Append one new data row (row 49) to Sheet1, mirroring the existing
# daily log rows: date (text), weekday (text), hour, ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the date-like string to stay text (matching
# the existing rows) instead of being auto-converted to a date serial.
$ws.Range("A49").Value = "'2025/10/02"
# Re-apply the plain/default style used by the rest of the data rows so
# the quote-prefix formatting doesn't leave behind an extra cell style.
$ws.Range("A49").Style = $ws.Range("A48").Style

$ws.Range("B49").Value = "木"
$ws.Range("C49").Value = 11
$ws.Range("D49").Value = 3
